$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '33.732.13'
$ws.Range('E2').Value = '  +8.73%  '
$ws.Range('D3').Value = '1.770.24'
$ws.Range('E3').Value = '  +4.55%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.57'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +1.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.552'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +3.45%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '30.44'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +3.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '46.34'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('E10').Value = '  +3.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0658'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +2.66%  '
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').Value = '2.023.52'
$ws.Range('E13').Value = '  +4.50%  '
$ws.Range('D14').Value = '1.770.46'
$ws.Range('E14').Value = '  +4.88%  '
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('D16').Value = '33.673.46'
$ws.Range('E16').Value = '  +8.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '9.96'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -2.43%  '
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '68.24'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '250.74'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +1.30%  '
$ws.Range('E21').Value = '  +1.96%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.22'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +1.14%  '
$ws.Range('E24').Value = '  -2.70%  '
$ws.Range('E25').Value = '  -1.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '158.22'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.43'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +3.13%  '
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.90'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +2.59%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  +5.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0511'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('E33').Value = '  +2.86%  '
$ws.Range('E34').Value = '  +4.92%  '
$ws.Range('D35').Value = '1.480.73'
$ws.Range('E35').Value = '  -2.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.79'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +3.34%  '
$ws.Range('E37').Value = '  +2.78%  '
$ws.Range('E38').Value = '  +2.25%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0184'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '83.01'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.35'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.67'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.882'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +4.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.07'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0510'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('E46').Value = '  +3.93%  '
$ws.Range('D47').Value = '1.921.42'
$ws.Range('E47').Value = '  +5.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.74'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +3.00%  '
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.75'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +14.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.34'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -2.57%  '
